# Adds the new "Probing" Capability and its child Strategic/Tactical
# Objectives to the malware_capabilities mind-map worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("malware_capabilities")

# Row 206: Capability "Probing"
$ws.Cells.Item(206, 1).Value = "Probing"
$ws.Cells.Item(206, 4).Value = "The 'probing' Capability indicates that the malware instance is able to probe its host system or network environment; most often this is done to support other Capabilities and their Objectives."

# Row 207: Strategic Objective "Probe Network Environment"
$ws.Cells.Item(207, 2).Value = "Probe Network Environment"
$ws.Cells.Item(207, 4).Value = "The 'probe network environment' value indicates that the malware instance is able to probe the properties of its network environment, e.g. to determine whether it funnels traffic through a proxy."

# Row 208: Tactical Objective "Check for Internet Connectivity"
$ws.Cells.Item(208, 3).Value = "Check for Internet Connectivity"
$ws.Cells.Item(208, 4).Value = "The 'check for internet connectivity' value indicates that the malware instance is able to check whether the network environment in which it executes is connected to the internet."

# Row 209: Tactical Objective "Check for Firewall"
$ws.Cells.Item(209, 3).Value = "Check for Firewall"
$ws.Cells.Item(209, 4).Value = "The 'check for firewall' value indicates that the malware instance is able to check whether the network environment in which it executes contains a hardware or software firewall."

# Row 210: Tactical Objective "Check for Proxy"
$ws.Cells.Item(210, 3).Value = "Check for Proxy"
$ws.Cells.Item(210, 4).Value = "The 'check for proxy' value indicates that the malware instance is able to check whether the network environment in which it executes contains a hardware or software proxy."

# Row 211: Tactical Objective "Map Local Network"
$ws.Cells.Item(211, 3).Value = "Map Local Network"
$ws.Cells.Item(211, 4).Value = "The 'map local network' value indicates that the malware instance is able to map the layout of the local network environment in which it executes."

# Row 212: Strategic Objective "Probe Host Configuration"
$ws.Cells.Item(212, 2).Value = "Probe Host Configuration"
$ws.Cells.Item(212, 4).Value = "The 'probe host configuration' value indicates that the malware instance is able to probe the configuration of the host system on which it executes."

# Row 213: Tactical Objective "Check Language"
$ws.Cells.Item(213, 3).Value = "Check Language"
$ws.Cells.Item(213, 4).Value = "The 'check language' value indicates that the malware instance is able to check the language of the host system on which it executes."

# Mirror the formatting of the existing rows (column A/B/C use style index 7 /
# "indent" look with border+left/top align, column D wraps text, E/F keep the
# bordered-blank look) by copying the row immediately above into the new rows,
# then overwriting the values above. Use the same left-border formats as the
# analogous "Capability" (row 202), "Strategic Objective" (row 200), and
# "Tactical Objective" (row 201) rows above it.
$wb.Worksheets.Item("malware_capabilities").Range("A202:F202").Copy() | Out-Null
$ws.Range("A206:F206").PasteSpecial(-4122) | Out-Null

$ws.Range("A201:F201").Copy() | Out-Null
$ws.Range("A207:F207").PasteSpecial(-4122) | Out-Null
$ws.Range("A207:F207").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Re-set values (copy/paste above only carried formats, values already set).
$ws.Cells.Item(206, 1).Value = "Probing"
$ws.Cells.Item(206, 4).Value = "The 'probing' Capability indicates that the malware instance is able to probe its host system or network environment; most often this is done to support other Capabilities and their Objectives."

# Scroll/selection state, matching the final saved view.
$ws.Application.Goto($ws.Range("A58"), $false)
$ws.Range("C5").Select() | Out-Null
